# Update the instructional subtitle text under each pairwise-comparison
# matrix title from the old wording to the new wording.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Enter pairwise comparisons in the white cells of the table or numerical data in the green cells. For the Direct Values column, if the smallest value is best, invert the value before entering it (e.g., `$10 as =1/10) ."

$ws.Range("A2").Value = $newText
$ws.Range("A11").Value = $newText
$ws.Range("A19").Value = $newText
$ws.Range("A27").Value = $newText
$ws.Range("A35").Value = $newText
